$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.503.92"
$ws.Range("E2").Value = "  -1.74%  "

$ws.Range("D3").Value = "1.669.04"
$ws.Range("E3").Value = "  -1.99%  "

$ws.Range("D4").Value = "'1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").Value = "'313.29"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.07%  "

$ws.Range("E6").Value = "  -0.03%  "

$ws.Range("D7").Value = "'0.3896"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.33%  "

$ws.Range("D8").Value = "'0.3925"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.33%  "

$ws.Range("D9").Value = "'1.003"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.02%  "

$ws.Range("D10").Value = "'51.63"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.82%  "

$ws.Range("D11").Value = "'1.403"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.75%  "

$ws.Range("D12").Value = "'0.08612"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.35%  "

$ws.Range("D13").Value = "'24.94"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.04%  "

$ws.Range("D14").Value = "'7.268"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.76%  "

$ws.Range("E15").Value = "  -3.07%  "

$ws.Range("D16").Value = "'7.707"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.13%  "

$ws.Range("D17").Value = "1.679.44"
$ws.Range("E17").Value = "  -2.24%  "

$ws.Range("D18").Value = "'93.34"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.08%  "

$ws.Range("D19").Value = "'0.07045"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.57%  "

$ws.Range("D20").Value = "'20.37"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.72%  "

$ws.Range("D21").Value = "'7.038"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.87%  "

$ws.Range("D22").Value = "'1.006"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.22%  "

$ws.Range("D23").Value = "'13.91"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.11%  "

$ws.Range("D24").Value = "24.503.50"
$ws.Range("E24").Value = "  -1.65%  "

$ws.Range("D25").Value = "'2.375"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.84%  "

$ws.Range("D26").Value = "'23.15"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.42%  "

$ws.Range("D27").Value = "'2.723"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.72%  "

$ws.Range("D28").Value = "'161.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.89%  "

$ws.Range("D29").Value = "'5.760"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -11.46%  "

$ws.Range("D30").Value = "'147.32"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.45%  "

$ws.Range("D31").Value = "'8.284"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.67%  "

$ws.Range("D32").Value = "'2.508"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +12.36%  "

$ws.Range("D33").Value = "1.860.57"
$ws.Range("E33").Value = "  -2.31%  "

$ws.Range("D34").Value = "'0.08291"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.45%  "

$ws.Range("B35").Value = "InternetComputer(DFINITY)"
$ws.Range("C35").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D35").Value = "'6.966"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.03%  "

$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").Value = "'0.03015"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.81%  "

$ws.Range("D37").Value = "'0.2788"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.79%  "

$ws.Range("D38").Value = "'0.9733"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.13%  "

$ws.Range("D39").Value = "'0.09434"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.88%  "

$ws.Range("D40").Value = "'1.525"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.28%  "

$ws.Range("D41").Value = "'10.25"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.53%  "

$ws.Range("D42").Value = "'0.7840"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.43%  "

$ws.Range("D43").Value = "'13.45"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.33%  "

$ws.Range("D44").Value = "'16.18"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -7.81%  "

$ws.Range("D45").Value = "'0.7067"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.87%  "

$ws.Range("D46").Value = "'2.539"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.49%  "

$ws.Range("D47").Value = "'4.172"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.82%  "

$ws.Range("E48").Value = "  +0.08%  "

$ws.Range("D49").Value = "'0.08575"
$ws.Range("D49").Style = "Normal"

$ws.Range("D50").Value = "'1.320"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.91%  "

$ws.Range("D51").Value = "'136.99"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.61%  "
